$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$np = $s.NotesPage
Write-Output $np
$props = $np | Get-Member
Write-Output $props
